$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the column D price cells we touch keep their original Text
# storage type (the values look numeric, so Excel would otherwise
# auto-convert them to Number on assignment).
$priceRows = @(2,3,4,5,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21,23,24,27,40,41,42,43,44,45,46,47,48,49,50,51)
foreach ($r in $priceRows) {
    $ws.Range("D$r").NumberFormat = "@"
}

$ws.Range("D2").Value = "259.01"
$ws.Range("D3").Value = "22.75"
$ws.Range("D4").Value = "6.172"
$ws.Range("D5").Value = "0.06088"
$ws.Range("D6").Value = "6.725"
$ws.Range("D7").Value = "3.476"
$ws.Range("D8").Value = "1.358"
$ws.Range("D9").Value = "0.7997"
$ws.Range("B10").Value = "One"
$ws.Range("C10").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D10").Value = "0.01326"
$ws.Range("E10").Value = "9OneONE"
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").Value = "0.1590"
$ws.Range("E11").Value = "10WazirXWRX"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").Value = "0.08051"
$ws.Range("E12").Value = "11MandalaExchangeTokenMDX"
$ws.Range("B13").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C13").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D13").Value = "0.03330"
$ws.Range("E13").Value = "12LiechtensteinCryptoassetsExchangeLCX"
$ws.Range("B14").Value = "BitrueCoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D14").Value = "0.03048"
$ws.Range("E14").Value = "13BitrueCoinBTR"
$ws.Range("B15").Value = "BitMartToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D15").Value = "0.09291"
$ws.Range("E15").Value = "14BitMartTokenBMX"
$ws.Range("B16").Value = "MCDex"
$ws.Range("C16").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D16").Value = "3.916"
$ws.Range("E16").Value = "15MCDexMCB"
$ws.Range("B17").Value = "BitForexToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D17").Value = "0.001700"
$ws.Range("E17").Value = "16BitForexTokenBF"
$ws.Range("B18").Value = "CoinExToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D18").Value = "0.04829"
$ws.Range("E18").Value = "17CoinExTokenCET"
$ws.Range("D19").Value = "0.006187"
$ws.Range("D20").Value = "0.001101"
$ws.Range("D21").Value = "0.003383"
$ws.Range("D23").Value = "3.693"
$ws.Range("D24").Value = "2.262"
$ws.Range("D27").Value = "0.0003020"
$ws.Range("D40").Value = "0.04593"
$ws.Range("D41").Value = "0.007171"
$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D42").Value = "0.003906"
$ws.Range("E42").Value = "41CEJICEJI"
$ws.Range("B43").Value = "BKEXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D43").Value = "0.1113"
$ws.Range("E43").Value = "42BKEXTokenBKK"
$ws.Range("D44").Value = "0.01062"
$ws.Range("D45").Value = "0.002973"
$ws.Range("D46").Value = "0.00005925"
$ws.Range("D47").Value = "0.00000000751"
$ws.Range("D48").Value = "0.7512"
$ws.Range("D49").Value = "0.06584"
$ws.Range("E49").Value = "48BOLOBOLOWorstin24h"
$ws.Range("D50").Value = "0.00001502"
$ws.Range("D51").Value = "0.01012"
